$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column in H, matching the formatting of the existing header row (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
